$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "20.365.75"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -6.63%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.439.10"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -6.62%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("E5").Value = "  -0.26%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "277.44"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3731"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -4.65%  "
$ws.Range("E8").Value = "  -3.20%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "40.62"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -5.74%  "
$ws.Range("E10").Value = "  -4.54%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.06600"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -7.73%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.18%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.374"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -4.40%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "17.32"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -6.60%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.149"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -7.19%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.436.44"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -7.15%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001011"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -7.89%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "76.61"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -8.06%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.05864"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -10.64%  "
$ws.Range("E20").Value = "  -0.20%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.735"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -6.55%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "14.38"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -5.66%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.99"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.328"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.64%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "20.358.11"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -6.71%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.250"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -5.26%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "142.82"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.27%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "17.06"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -7.15%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.603.01"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -6.84%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "110.19"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -5.69%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.974"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -18.07%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.9245"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -4.32%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.495"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -5.86%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.07720"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -5.99%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "8.346"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -7.35%  "
$ws.Range("E36").Value = "  +3.70%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.05731"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -5.89%  "
$ws.Range("E38").Value = "  -0.25%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "4.752"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -6.58%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.136"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.65%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.1928"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -5.10%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.02030"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -8.88%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.342"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -11.39%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.591"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -4.08%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.5345"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -6.79%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "12.08"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -5.67%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.5177"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -6.32%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "112.23"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.63%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.788"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.50%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.058"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -6.12%  "
$ws.Range("E51").Value = "  -0.30%  "
